$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G7").Value = 1.73
$ws.Range("H7").Value = 3.4
$ws.Range("I7").Value = 5
$ws.Range("M7").Value = 1.1
$ws.Range("N7").Value = 7
$ws.Range("O7").Value = 1.44
$ws.Range("P7").Value = 2.63
$ws.Range("AC7").Value = 7
